# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.998.75"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "3.036.95"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'595.18"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'154.15"
$ws.Range("E6").Value = "  +7.31%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.032.54"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'6.87"
$ws.Range("E10").Value = "  +13.59%  "
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").Value = "'0.0000235"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "'35.87"
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("D15").Value = "'0.125"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "3.538.88"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "62.926.98"
$ws.Range("D19").Value = "3.037.92"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "'454.44"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'14.30"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").Value = "'83.12"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'11.33"
$ws.Range("E25").Value = "  +7.31%  "
$ws.Range("E26").Value = "  +4.53%  "
$ws.Range("E27").Value = "  +4.18%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'7.52"
$ws.Range("E29").Value = "  +5.86%  "
$ws.Range("E30").Value = "  +10.08%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'27.63"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "0.0₃0861"
$ws.Range("E35").Value = "  +4.82%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("E37").Value = "  +3.40%  "
$ws.Range("D38").Value = "'3.20"
$ws.Range("E38").Value = "  +11.30%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'2.12"
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.131"
$ws.Range("E40").Value = "  +6.61%  "
$ws.Range("D41").Value = "'50.35"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'0.307"
$ws.Range("E43").Value = "  +14.05%  "
$ws.Range("D44").Value = "'44.09"
$ws.Range("E44").Value = "  +11.19%  "
$ws.Range("D45").Value = "'392.23"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D47").Value = "2.724.40"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "'132.82"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("E50").Value = "  +7.11%  "
$ws.Range("D51").Value = "'24.67"
$ws.Range("E51").Value = "  +5.32%  "
